# Apply rounding function for coefficients (xls)
#
# Source data (rows 34-36) and the deflection-constant cell (E30) were
# updated with new measurements, and the coefficient formulas in column C
# (C7:C14) were changed to wrap the existing calculation in ROUND(...,2),
# with a "0.000" number format applied so the rounded coefficient is shown
# with three decimals.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Update the raw load-cell source table (rows 34-36) ---
$ws.Range("D34").Value = 23713
$ws.Range("E34").Value = 73305

$ws.Range("D35").Value = -171611
$ws.Range("E35").Value = -122443

$ws.Range("D36").Value = 0
$ws.Range("E36").Value = 0

# --- Update the deflection constant used by the coefficient formula ---
$ws.Range("E30").Value = 1.219

# --- Rewrite the column-C coefficient formulas to round to 2 decimals ---
$ws.Range("C7").Formula  = "=ROUND((E34-D34)/(9.81*`$E`$30),2)"
$ws.Range("C8").Formula  = "=ROUND((E35-D35)/(9.81*`$E`$30),2)"
$ws.Range("C9").Formula  = "=ROUND((E36-D36)/(9.81*`$E`$30),2)"
$ws.Range("C10").Formula = "=ROUND((E37-D37)/(9.81*`$E`$30),2)"
$ws.Range("C11").Formula = "=ROUND((E38-D38)/(9.81*`$E`$30),2)"
$ws.Range("C12").Formula = "=ROUND((E39-D39)/(9.81*`$E`$30),2)"
$ws.Range("C13").Formula = "=ROUND((E40-D40)/(9.81*`$E`$30),2)"
$ws.Range("C14").Formula = "=ROUND((E41-D41)/(9.81*`$E`$30),2)"

# New style for the rounded coefficients: 3-decimal number format, same
# green fill / border / centered alignment the cells already had.
$ws.Range("C7:C14").NumberFormat = "0.000"

# --- Selection moved by the author while editing ---
$ws.Range("J17").Select()
